# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.
#
# For each affected sheet:
#   1. Fix the E1 header label (was a leftover numeric value 687.0428...,
#      should be the text label for the last year/period column).
#   2. Remove the trailing "Total" row.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellAddr, $text) {
    # Assign as a text formula first, then paste-special as values only so
    # that the result is stored as a genuine text value (not a number and
    # not a formula) while keeping the cell's existing style untouched.
    $range = $ws.Range($cellAddr)
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# Sheet 1: "Potencia Acumulada - SIN (MW)" -> E1 = "2050", remove row 13 (Total)
$ws1 = $wb.Worksheets.Item("Potencia Acumulada - SIN (MW)")
Set-TextValue $ws1 "E1" "2050"
$ws1.Rows.Item(13).Delete()

# Sheet 2: "Geracao Periodo Medio (MWMed)" -> E1 = "2050", remove row 13 (Total)
$ws2 = $wb.Worksheets.Item("Geracao Periodo Medio (MWMed)")
Set-TextValue $ws2 "E1" "2050"
$ws2.Rows.Item(13).Delete()

# Sheet 3: "Atendimento a Ponta(MW)" -> E1 = "2050", remove row 13 (Total)
$ws3 = $wb.Worksheets.Item("Atendimento a Ponta(MW)")
Set-TextValue $ws3 "E1" "2050"
$ws3.Rows.Item(13).Delete()

# Sheet 4: "Potencia Incremental - SIN(MW)" -> E1 = "2041-2050", remove row 13 (Total)
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextValue $ws4 "E1" "2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5: "Emissoes Totais (MtCO2eq)" -> unchanged

# Sheet 6: "Custo Total (bilhões de R$)" -> remove row 4 (Total)
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
